# Apply "final output data" update: append two new data rows (12 & 13)
# to the "Execution Times" sheet, matching the new rows added to the
# source XLSX in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Execution Times")

# Row 12 - Graph DFS (col C/D) and A* (col K/L) timings
$ws.Range("C12").Value = [double]"0.006205"
$ws.Range("D12").Value = [double]"0.0010974"
$ws.Range("K12").Value = [double]"2.724E-4"
$ws.Range("L12").Value = [double]"4.16E-4"

# Row 13 - Graph DFS (col C/D) timings
$ws.Range("C13").Value = [double]"0.0022576"
$ws.Range("D13").Value = [double]"0.005929"
